$p = $ppt.ActivePresentation

# Slide 3 (pptx slide index 3) -> "Slide 2: History of Domestic Cats"
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 2: History of Domestic Cats"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats have been domesticated for thousands of years and have played various roles in human society. They were first domesticated in the Near East around 7500 BC and have since spread to all corners of the globe. Today, cats are one of the most popular pets in the world."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://www.nationalgeographic.com/animals/mammals/d/domestic-cat/"

# Slide 4 (pptx slide index 4) -> "Slide 3: Physical Characteristics of Cats"
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 3: Physical Characteristics of Cats"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats have a flexible body, sharp retractable claws, strong jaws, and excellent vision and hearing. They are known for their retractable claws that they use for hunting and climbing. Cats also have a keen sense of smell and can see well in low light."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://www.thesprucepets.com/cat-breeds-4176798"

# Slide 5 (pptx slide index 5) -> "Slide 4: Cat Behavior and Communication"
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 4: Cat Behavior and Communication"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats communicate through a combination of vocalizations, body language, and scent marking. They can purr when content, meow when hungry or seeking attention, and hiss or growl when threatened. Understanding cat behavior is important for creating a strong bond with your feline friend."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://icatcare.org/advice/cat-behaviour/"

# Slide 6 (pptx slide index 6) -> "Slide 5: Health and Care of Cats"
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 5: Health and Care of Cats"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats require regular grooming, vaccinations, and veterinary check-ups to ensure they stay healthy. It's important to provide them with a balanced diet, fresh water, and a safe environment. Proper care can help prevent common health issues in cats."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://pets.webmd.com/cats/guide/caring-for-a-cat"

# Slide 7 (pptx slide index 7) -> "Slide 6: Common Breeds of Cats"
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 6: Common Breeds of Cats"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "There are over 70 recognized cat breeds, each with its own unique characteristics and personality traits. Some popular breeds include the Siamese, Persian, Maine Coon, and Bengal. Different breeds have different grooming and care requirements."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://cattime.com/cat-breeds"

# Slide 8 (pptx slide index 8) -> "Slide 7: Famous Cats in History and Pop Culture"
$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 7: Famous Cats in History and Pop Culture"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats have played important roles in history and pop culture. From ancient Egyptian worship of cats to famous fictional felines like Garfield and Hello Kitty, cats have captured the hearts of people around the world."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://www.history.com/news/a-history-of-house-cats"

# Slide 9 (pptx slide index 9) -> "Slide 8: Cat Myths and Superstitions"
$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 8: Cat Myths and Superstitions"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats have been associated with various myths and superstitions throughout history. In many cultures, black cats are considered bad luck, while in others, they are seen as symbols of good fortune. Understanding these myths can help dispel common misconceptions about cats."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://www.livescience.com/35041-black-cat-taboo-life-balance-spooky-110713.html"

# Slide 10 (pptx slide index 10) -> "Slide 9: Cats in Art and Literature"
$s = $p.Slides.Item(10)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Slide 9: Cats in Art and Literature"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats have been a popular subject in art and literature for centuries. From ancient Egyptian hieroglyphics to modern-day internet memes, cats continue to inspire creativity and imagination. Many famous artists and writers have featured cats in their work."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://www.artsy.net/article/artsy-editorial-art-historian-cats-famously-feline-history-art"

# Slide 11 (pptx slide index 11) -> "Slide 10: Conclusion" (title unchanged)
$s = $p.Slides.Item(11)
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Cats are fascinating creatures with a long history of companionship with humans. Whether as pets, symbols of luck, or characters in stories, cats continue to hold a special place in our hearts. By understanding their behavior, health needs, and unique traits, we can build strong bonds with our feline friends."
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1,1).Text = "- https://www.animalwised.com/history-and-origin-of-the-domestic-cat-3609.html"
